$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, derived from the authoritative diff.
$updates = [ordered]@{
    'D2' = '27.028.65'
    'E2' = '  +0.54%  '
    'D3' = '1.888.34'
    'E3' = '  +1.44%  '
    'E4' = '  +0.10%  '
    'D5' = '306.14'
    'E5' = '  +0.45%  '
    'D6' = '1.001'
    'D7' = '0.5186'
    'E7' = '  +2.89%  '
    'D8' = '0.3754'
    'E8' = '  +2.99%  '
    'D9' = '0.07207'
    'E9' = '  +0.57%  '
    'E10' = '  +2.17%  '
    'E11' = '  +1.04%  '
    'D12' = '0.07629'
    'E12' = '  +1.72%  '
    'D13' = '1.872.24'
    'E13' = '  +0.30%  '
    'D14' = '94.50'
    'E14' = '  -0.46%  '
    'D15' = '5.235'
    'E15' = '  +0.11%  '
    'D16' = '1.001'
    'E16' = '  +0.08%  '
    'D17' = '0.000008509'
    'E17' = '  -0.02%  '
    'D18' = '14.44'
    'E18' = '  +1.63%  '
    'D19' = '1.000'
    'E19' = '  +0.01%  '
    'D20' = '27.078.68'
    'E20' = '  +0.52%  '
    'D22' = '2.116.49'
    'E22' = '  +0.15%  '
    'D23' = '10.60'
    'E23' = '  +2.08%  '
    'D24' = '6.386'
    'E24' = '  -0.25%  '
    'E25' = '  +10.21%  '
    'D26' = '145.65'
    'E26' = '  -1.45%  '
    'D27' = '1.741'
    'E27' = '  -2.32%  '
    'E28' = '  +1.15%  '
    'D29' = '114.09'
    'E29' = '  +0.80%  '
    'D30' = '4.911'
    'E30' = '  +5.26%  '
    'D31' = '4.791'
    'E31' = '  +1.80%  '
    'D32' = '0.09195'
    'E32' = '  -0.25%  '
    'D33' = '0.05036'
    'E33' = '  -2.06%  '
    'E34' = '  +7.55%  '
    'D35' = '0.7667'
    'E35' = '  +2.26%  '
    'D36' = '2.961'
    'E36' = '  +0.38%  '
    'D37' = '3.276'
    'E37' = '  +0.45%  '
    'D38' = '2.599'
    'E38' = '  +0.53%  '
    'D39' = '0.5597'
    'E39' = '  +0.51%  '
    'D40' = '0.01987'
    'E40' = '  -0.68%  '
    'D41' = '1.070'
    'E41' = '  +0.09%  '
    'D42' = '9.013'
    'E42' = '  +4.95%  '
    'D43' = '6.628'
    'E43' = '  +1.10%  '
    'D44' = '119.02'
    'E44' = '  +2.24%  '
    'D45' = '0.1506'
    'E45' = '  +2.38%  '
    'D46' = '0.4823'
    'B47' = 'PaxDollar'
    'C47' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'D47' = '1.000'
    'E47' = '  +0.09%  '
    'B48' = 'EnergySwap'
    'C48' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D48' = '10.13'
    'E48' = '  +1.03%  '
    'D49' = '1.597'
    'E49' = '  +2.32%  '
    'D50' = '37.74'
    'E50' = '  +2.89%  '
    'E51' = '  +1.52%  '
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "1.000", "0.07207",
    # thousand-dot prices like "27.028.65") are not reinterpreted as numbers
    # and keep their exact original formatting/precision.
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
}
